$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused duplicate chart defined names (v1.2 / v1.3)
$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

# New headers for the "increase" stats block
$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

# Mean increase vs. a reference mean, and median increase vs. a reference median
$ws.Range("D19").Formula = "= ((E3 / 114.202998) * 100) - 100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").Formula = "= ((E10 / 113.658804) * 100) - 100"
$ws.Range("F19").ClearFormats()

[void]$ws.Range("E20").Select()
